$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-07 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-08 Friday", 2)

$d.Content.Find.Execute("715×4=", $true, $false, $false, $false, $false, $true, 1, $false, "950×2=", 2)
$d.Content.Find.Execute("104×8=", $true, $false, $false, $false, $false, $true, 1, $false, "387×5=", 2)
$d.Content.Find.Execute("571×5=", $true, $false, $false, $false, $false, $true, 1, $false, "115×4=", 2)
$d.Content.Find.Execute("878×9=", $true, $false, $false, $false, $false, $true, 1, $false, "555×8=", 2)
$d.Content.Find.Execute("370×2=", $true, $false, $false, $false, $false, $true, 1, $false, "887×9=", 2)

$d.Content.Find.Execute("112×2=", $true, $false, $false, $false, $false, $true, 1, $false, "466×9=", 2)
$d.Content.Find.Execute("603×9=", $true, $false, $false, $false, $false, $true, 1, $false, "669×8=", 2)
$d.Content.Find.Execute("666×5=", $true, $false, $false, $false, $false, $true, 1, $false, "261×9=", 2)
$d.Content.Find.Execute("945×7=", $true, $false, $false, $false, $false, $true, 1, $false, "356×3=", 2)
$d.Content.Find.Execute("674×6=", $true, $false, $false, $false, $false, $true, 1, $false, "753×2=", 2)

$d.Content.Find.Execute("518×6=", $true, $false, $false, $false, $false, $true, 1, $false, "236×4=", 2)
$d.Content.Find.Execute("934×3=", $true, $false, $false, $false, $false, $true, 1, $false, "356×4=", 2)
$d.Content.Find.Execute("565×6=", $true, $false, $false, $false, $false, $true, 1, $false, "910×2=", 2)
$d.Content.Find.Execute("916×2=", $true, $false, $false, $false, $false, $true, 1, $false, "302×9=", 2)
$d.Content.Find.Execute("500×7=", $true, $false, $false, $false, $false, $true, 1, $false, "318×5=", 2)

$d.Content.Find.Execute("476×3=", $true, $false, $false, $false, $false, $true, 1, $false, "853×8=", 2)
$d.Content.Find.Execute("860×5=", $true, $false, $false, $false, $false, $true, 1, $false, "580×3=", 2)
$d.Content.Find.Execute("612×8=", $true, $false, $false, $false, $false, $true, 1, $false, "963×4=", 2)
$d.Content.Find.Execute("961×4=", $true, $false, $false, $false, $false, $true, 1, $false, "731×9=", 2)
$d.Content.Find.Execute("183×5=", $true, $false, $false, $false, $false, $true, 1, $false, "779×6=", 2)

$d.Content.Find.Execute("711×3=", $true, $false, $false, $false, $false, $true, 1, $false, "416×8=", 2)
$d.Content.Find.Execute("349×9=", $true, $false, $false, $false, $false, $true, 1, $false, "169×3=", 2)
$d.Content.Find.Execute("172×2=", $true, $false, $false, $false, $false, $true, 1, $false, "982×2=", 2)
$d.Content.Find.Execute("143×8=", $true, $false, $false, $false, $false, $true, 1, $false, "554×7=", 2)
$d.Content.Find.Execute("335×9=", $true, $false, $false, $false, $false, $true, 1, $false, "402×4=", 2)
